# UserController und Registration View erstellen
# - Mark the "Klasse UserRepositoryDB erstellen" task as "done" (was "b")
# - Mark the new "Registration Methode im UserController erstellen" and
#   "Registration View erstellen" tasks as "done"
# - Move the active selection to the last-edited cell (C14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("C12").Value = "done"
$ws.Range("C13").Value = "done"
$ws.Range("C14").Value = "done"

[void]$ws.Range("C14").Select()
